# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet
#    and fill it with the per-fund holding detail for 2022-Q1.
# 2. Insert a new top data row into the "总计" sheet summarising the
#    2022-Q1 quarter (9 funds, 1.61 亿元) and bump the existing running
#    index column (A) down by one for every pre-existing row.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" worksheet (placed immediately before 总计)
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet variables in this host resolve by tab *position*, not
# stable identity. Inserting "2022-Q1" shifted "总计" one slot to the
# right, so the old `$total` handle now aliases the new sheet instead.
# Re-fetch "总计" by name before touching it again.
$total = $wb.Worksheets.Item("总计")

# Reuse the exact header/column formatting from an existing quarter sheet
# (bold + bordered style) so the new sheet matches the others exactly.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$data = @(
    @("010695","华夏磐益一年定期开放混合","18.02","82.41","3.14","0.5658",3),
    @("009837","华夏磐锐一年定期开放混合A","16.45","79.44","3.13","0.5149",6),
    @("360001","光大保德信量化股票","17.53","88.21","1.83","0.3208",7),
    @("008895","申万菱信量化对冲策略灵活配置混合","8.57","42.67","1.02","0.0874",9),
    @("014135","中欧金安量化混合A","9.43","67.44","0.61","0.0575",9),
    @("011231","光大保德信锦弘混合A","4.13","20.96","0.70","0.0289",10),
    @("009838","华夏磐锐一年定期开放混合C","0.44","79.44","3.13","0.0138",6),
    @("011232","光大保德信锦弘混合C","1.29","20.96","0.70","0.0090",10),
    @("014136","中欧金安量化混合C","1.28","67.44","0.61","0.0078",9)
)

# Force columns B:G as text (number-formats are set to "@" first, values
# written, then formats cleared) so fund codes keep leading zeros and the
# ratio/percentage columns stay literal text like the source workbook.
$newSheet.Range("B2:G10").NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $newSheet.Range("A" + $r).Value = $r - 2
    $newSheet.Range("B" + $r).Value = $row[0]
    $newSheet.Range("C" + $r).Value = $row[1]
    $newSheet.Range("D" + $r).Value = $row[2]
    $newSheet.Range("E" + $r).Value = $row[3]
    $newSheet.Range("F" + $r).Value = $row[4]
    $newSheet.Range("G" + $r).Value = $row[5]
    $newSheet.Range("H" + $r).Value = $row[6]
    $r = $r + 1
}

$newSheet.Range("B2:G10").ClearFormats()

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 1.61

# Bump the running index (column A) of every pre-existing row by one.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
